$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

$ws.Range("I3").Value = "get date"
$ws.Range("I4").Value = "get amount"
$ws.Range("I5").Value = "set balance"

$ws.Range("H3").Value = "ar"
$ws.Range("H4").Value = "ar"
$ws.Range("H4").NumberFormat = "#,##0"

$ws.Range("H6").Value = "pa"
$ws.Range("H7").Value = "pa"

$ws.Range("I6").Value = "get date"
$ws.Range("I7").Value = "get amount"
$ws.Range("I8").Value = "set balance"

$ws.Range("K4").Value = 5000
$ws.Range("K5").Value = 5000
$ws.Range("K7").Value = 1000
$ws.Range("K8").Formula = "=+K5-K7"

$ws.Range("A3").Select() | Out-Null
